$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 581.3333
$ws.Range("I12").Value = 446.66666
$ws.Range("J12").Value = 648.6667
$ws.Range("K12").Value = 446.66666
$ws.Range("L12").Value = 648.6667
$ws.Range("M12").Value = -276.66666
$ws.Range("N12").Value = -988.6667
$ws.Range("H32").Value = 4728.2
$ws.Range("J32").Value = 5554.5
$ws.Range("L32").Value = 5554.5
$ws.Range("N32").Value = -6206.5
$ws.Range("H113").Value = 7582.8335
$ws.Range("I113").Value = 6874.5
$ws.Range("J113").Value = 8999.5
$ws.Range("K113").Value = 6874.5
$ws.Range("L113").Value = 8999.5
$ws.Range("M113").Value = -3620.5
$ws.Range("N113").Value = -15507.5
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 2628.9778
$ws.Range("I132").Value = 2671.9023
$ws.Range("J132").Value = 2189
$ws.Range("K132").Value = 8015.706900000001
$ws.Range("L132").Value = 6567
$ws.Range("M132").Value = -5485.706900000001
$ws.Range("N132").Value = -11627
$ws.Range("H137").Value = 1295.8572
$ws.Range("I137").Value = 799
$ws.Range("J137").Value = 1378.6666
$ws.Range("K137").Value = 2397
$ws.Range("L137").Value = 4135.9998
$ws.Range("M137").Value = 153
$ws.Range("N137").Value = -9235.9998

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1892.9661
$ws.Range("I32").Value = 1030.1964
$ws.Range("J32").Value = 17998
$ws.Range("K32").Value = 1030.1964
$ws.Range("L32").Value = 17998
$ws.Range("M32").Value = -743.1964
$ws.Range("N32").Value = -18572
$ws.Range("H74").Value = 1966.5
$ws.Range("I74").Value = 1837.375
$ws.Range("J74").Value = 2999.5
$ws.Range("K74").Value = 1837.375
$ws.Range("L74").Value = 2999.5
$ws.Range("M74").Value = -963.375
$ws.Range("N74").Value = -4747.5
$ws.Range("H77").Value = 1966.5
$ws.Range("I77").Value = 1837.375
$ws.Range("J77").Value = 2999.5
$ws.Range("K77").Value = 9186.875
$ws.Range("L77").Value = 14997.5
$ws.Range("M77").Value = -4818.875
$ws.Range("N77").Value = -23733.5
$ws.Range("H88").Value = 1599.2693
$ws.Range("I88").Value = 1268
$ws.Range("J88").Value = 1774.6471
$ws.Range("K88").Value = 1268
$ws.Range("L88").Value = 1774.6471
$ws.Range("M88").Value = -862
$ws.Range("N88").Value = -2586.6471
$ws.Range("H91").Value = 1599.2693
$ws.Range("I91").Value = 1268
$ws.Range("J91").Value = 1774.6471
$ws.Range("K91").Value = 1268
$ws.Range("L91").Value = 1774.6471
$ws.Range("M91").Value = 136
$ws.Range("N91").Value = -4582.6471
$ws.Range("H97").Value = 485.68
$ws.Range("I97").Value = 423.56522
$ws.Range("K97").Value = 423.56522
$ws.Range("M97").Value = 72.43477999999999
$ws.Range("H110").Value = 2766.6667
$ws.Range("I110").Value = 2750
$ws.Range("J110").Value = 2800
$ws.Range("K110").Value = 2750
$ws.Range("L110").Value = 2800
$ws.Range("M110").Value = -705
$ws.Range("N110").Value = -6890
$ws.Range("H122").Value = 5699.778
$ws.Range("I122").Value = 5666
$ws.Range("J122").Value = 5970
$ws.Range("K122").Value = 16998
$ws.Range("L122").Value = 17910
$ws.Range("M122").Value = -14548
$ws.Range("N122").Value = -22810
$ws.Range("H132").Value = 5807.087
$ws.Range("I132").Value = 5148.386
$ws.Range("J132").Value = 20298.5
$ws.Range("K132").Value = 15445.158
$ws.Range("L132").Value = 60895.5
$ws.Range("M132").Value = -12915.158
$ws.Range("N132").Value = -65955.5

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1823.0769
$ws.Range("I86").Value = 1913.8
$ws.Range("J86").Value = 1520.6666
$ws.Range("K86").Value = 1913.8
$ws.Range("L86").Value = 1520.6666
$ws.Range("M86").Value = -790.8
$ws.Range("N86").Value = -3766.6666
$ws.Range("H89").Value = 1823.0769
$ws.Range("I89").Value = 1913.8
$ws.Range("J89").Value = 1520.6666
$ws.Range("K89").Value = 9569
$ws.Range("L89").Value = 7603.333000000001
$ws.Range("M89").Value = -3953
$ws.Range("N89").Value = -18835.333

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5571.425
$ws.Range("I31").Value = 2274.25
$ws.Range("J31").Value = 6815.6416
$ws.Range("K31").Value = 2274.25
$ws.Range("L31").Value = 6815.6416
$ws.Range("M31").Value = -1979.25
$ws.Range("N31").Value = -7405.6416
$ws.Range("H34").Value = 5571.425
$ws.Range("I34").Value = 2274.25
$ws.Range("J34").Value = 6815.6416
$ws.Range("K34").Value = 2274.25
$ws.Range("L34").Value = 6815.6416
$ws.Range("M34").Value = -2072.25
$ws.Range("N34").Value = -7219.6416
$ws.Range("H58").Value = 4421.381
$ws.Range("I58").Value = 2963.8667
$ws.Range("J58").Value = 8065.1665
$ws.Range("K58").Value = 2963.8667
$ws.Range("L58").Value = 8065.1665
$ws.Range("M58").Value = -2760.8667
$ws.Range("N58").Value = -8471.166499999999
$ws.Range("H122").Value = 998.3333
$ws.Range("I122").Value = 706
$ws.Range("J122").Value = 1144.5
$ws.Range("K122").Value = 2118
$ws.Range("L122").Value = 3433.5
$ws.Range("M122").Value = 332
$ws.Range("N122").Value = -8333.5
$ws.Range("H132").Value = 1625.6765
$ws.Range("I132").Value = 1490.0344
$ws.Range("J132").Value = 2412.4
$ws.Range("K132").Value = 4470.1032
$ws.Range("L132").Value = 7237.200000000001
$ws.Range("M132").Value = -1940.1032
$ws.Range("N132").Value = -12297.2
$ws.Range("H136").Value = 4421.381
$ws.Range("I136").Value = 2963.8667
$ws.Range("J136").Value = 8065.1665
$ws.Range("K136").Value = 8891.6001
$ws.Range("L136").Value = 24195.4995
$ws.Range("M136").Value = -6341.6001
$ws.Range("N136").Value = -29295.4995

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 588313
$ws.Range("I2").Value = 869600.25
$ws.Range("J2").Value = 167
$ws.Range("K2").Value = 5217601.5
$ws.Range("L2").Value = 1002
$ws.Range("M2").Value = -5217488.5
$ws.Range("N2").Value = -1228
$ws.Range("H68").Value = 2761.3333
$ws.Range("J68").Value = 2820.9546
$ws.Range("L68").Value = 8462.863799999999
$ws.Range("N68").Value = -10084.8638
$ws.Range("H71").Value = 2761.3333
$ws.Range("J71").Value = 2820.9546
$ws.Range("L71").Value = 25388.5914
$ws.Range("N71").Value = -33500.5914
$ws.Range("H107").Value = 3105.7856
$ws.Range("J107").Value = 3498.4443
$ws.Range("L107").Value = 10495.3329
$ws.Range("N107").Value = -14335.3329

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 29201.285
$ws.Range("I70").Value = 61125.668
$ws.Range("K70").Value = 61125.668
$ws.Range("M70").Value = -60855.668
$ws.Range("H73").Value = 29201.285
$ws.Range("I73").Value = 61125.668
$ws.Range("K73").Value = 61125.668
$ws.Range("M73").Value = -60189.668
$ws.Range("H80").Value = 2676.923
$ws.Range("I80").Value = 1527.2727
$ws.Range("J80").Value = 9000
$ws.Range("K80").Value = 1527.2727
$ws.Range("L80").Value = 9000
$ws.Range("M80").Value = -529.2727
$ws.Range("N80").Value = -10996
$ws.Range("H83").Value = 2676.923
$ws.Range("I83").Value = 1527.2727
$ws.Range("J83").Value = 9000
$ws.Range("K83").Value = 7636.363499999999
$ws.Range("L83").Value = 45000
$ws.Range("M83").Value = -2644.363499999999
$ws.Range("N83").Value = -54984
$ws.Range("H122").Value = 3075.56
$ws.Range("I122").Value = 3127.111
$ws.Range("J122").Value = 2943
$ws.Range("K122").Value = 9381.332999999999
$ws.Range("L122").Value = 8829
$ws.Range("M122").Value = -6931.332999999999
$ws.Range("N122").Value = -13729
$ws.Range("H133").Value = 67250
$ws.Range("J133").Value = 67250
$ws.Range("L133").Value = 67250
$ws.Range("N133").Value = -77370

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 48639.668
$ws.Range("J43").Value = 22959.5
$ws.Range("L43").Value = 22959.5
$ws.Range("N43").Value = -23345.5
$ws.Range("H61").Value = 4107.84
$ws.Range("I61").Value = 4230.304
$ws.Range("J61").Value = 2699.5
$ws.Range("K61").Value = 4230.304
$ws.Range("L61").Value = 2699.5
$ws.Range("M61").Value = -4028.304
$ws.Range("N61").Value = -3103.5
$ws.Range("H113").Value = 4107.84
$ws.Range("I113").Value = 4230.304
$ws.Range("J113").Value = 2699.5
$ws.Range("K113").Value = 4230.304
$ws.Range("L113").Value = 2699.5
$ws.Range("M113").Value = -2060.304
$ws.Range("N113").Value = -7039.5
$ws.Range("H122").Value = 3498.8965
$ws.Range("I122").Value = 3638.64
$ws.Range("J122").Value = 2625.5
$ws.Range("K122").Value = 10915.92
$ws.Range("L122").Value = 7876.5
$ws.Range("M122").Value = -8465.92
$ws.Range("N122").Value = -12776.5
$ws.Range("H136").Value = 2288.0784
$ws.Range("I136").Value = 1891.6154
$ws.Range("J136").Value = 2700.4
$ws.Range("K136").Value = 5674.8462
$ws.Range("L136").Value = 8101.200000000001
$ws.Range("M136").Value = -3124.8462
$ws.Range("N136").Value = -13201.2

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 257.0476
$ws.Range("I107").Value = 274.94446
$ws.Range("J107").Value = 149.66667
$ws.Range("K107").Value = 824.83338
$ws.Range("L107").Value = 449.00001
$ws.Range("M107").Value = 1095.16662
$ws.Range("N107").Value = -4289.00001
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H136").Value = 3240.4285
$ws.Range("I136").Value = 3308.7646
$ws.Range("J136").Value = 2950
$ws.Range("K136").Value = 9926.293799999999
$ws.Range("L136").Value = 8850
$ws.Range("M136").Value = -7376.293799999999
$ws.Range("N136").Value = -13950
